# Add a "ManualCheck" column (C) to Sheet1, populate it for all existing
# data rows, and remove the trailing "Ctep" row (old row 86) that the
# commit migrated out of the hardcoded list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First drop the last row ("Ctep" / http://ctep-dev.wppro.lacounty.gov/)
# before writing column C so the remaining data rows keep their original
# row numbers (1-85).
$ws.Rows.Item(86).Delete()

# Header + per-row "Yes"/"No" values for the new ManualCheck column,
# in row order (row 1 is the header).
$manualCheckValues = @("ManualCheck","No","No","No","No","No","No","Yes","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","Yes","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","Yes","No","No","No","No","Yes","No","No","Yes","Yes","No","No","Yes","No","No","No","No","No","No","No","No","Yes","Yes","No","No","No","No","No","No","No","Yes","No","No")

for ($i = 0; $i -lt $manualCheckValues.Length; $i++) {
    $rowNum = $i + 1
    $ws.Cells.Item($rowNum, 3).Value = $manualCheckValues[$i]
}

Write-Output "ManualCheck column written for $($manualCheckValues.Length) rows; row 86 removed."
